$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header in B1
$ws.Range("B1").Value = "upper case"

# Fill in column B with upper-case values, and update column A values
$ws.Range("A2").Value = "mno"
$ws.Range("B2").Value = "MNO"

$ws.Range("A3").Value = "xyz"
$ws.Range("B3").Value = "XYZ"

$ws.Range("A4").Value = "wsg"
$ws.Range("B4").Value = "WSG"

$ws.Range("A5").Value = "lpo"
$ws.Range("B5").Value = "LPO"

$ws.Range("A6").Value = "mno"
$ws.Range("B6").Value = "MNO"

# Update selection to H7
$ws.Range("H7").Select()

$wb.Save()
